# Daily attendance processing - 2025-11-04 08:28:08
# Reverses the order of names/emails in the "Recorded By" column (G) for the
# specific session rows touched by today's processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRows = @(
    2,3,5,6,7,8,10,12,13,14,15,18,19,20,21,22,24,26,
    28,29,31,32,33,34,36,38,39,40,41,44,45,46,47,48,50,52,
    54,55,57,58,59,60,62,64,65,66,67,70,71,72,73,74,76,78,
    80,81,82,83,84,85,86,87,90,92,99,101,
    106,107,108,109,110,111,112,113,116,118,125,127,
    132,133,134,135,136,137,138,139,142,144,151,153
)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value()
    if ($current) {
        $parts = $current -split ', '
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value = ($reversed -join ', ')
        }
    }
}
